$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F1: year range header ---
$ws.Range("F1").Value = "2023-2024"

# --- Row 2: month headers ---
# Unmerge the existing month-header merges before rebuilding the new layout
$ws.Range("F2:J2").UnMerge()
$ws.Range("K2:N2").UnMerge()
$ws.Range("O2:P2").UnMerge()

# G2, L2 and P2 are brand-new merge anchors that previously sat inside a
# merged range (so they carry no explicit style) - copy the header style
# from F2 (which already kept its style as the old merge's anchor cell)
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("P2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F2").Value = "December"
$ws.Range("G2").Value = "January"
$ws.Range("L2").Value = "February"
$ws.Range("P2").Value = "March"

$ws.Range("G2:K2").Merge()
$ws.Range("L2:O2").Merge()

# --- Row 3: weekly date ranges ---
$ws.Range("F3").Value = "25/Dec - 31/Dec"
$ws.Range("G3").Value = "01/Jan - 07/Jan"
$ws.Range("H3").Value = "08/Jan - 14/Jan"
$ws.Range("I3").Value = "15/Jan - 21/Jan"
$ws.Range("J3").Value = "22/Jan - 28/Jan"
$ws.Range("K3").Value = "29/Jan - 04/Feb"
$ws.Range("L3").Value = "05/Feb - 11/Feb"
$ws.Range("M3").Value = "12/Feb - 18/Feb"
$ws.Range("N3").Value = "19/Feb - 25/Feb"
$ws.Range("O3").Value = "26/Feb - 03/Mar"
$ws.Range("P3").Value = "04/Mar - 10/Mar"

# --- Rows 4-13: Start Date / End Date columns ---
$ws.Range("D4").Value = "12/25"
$ws.Range("E4").Value = "12/31"

$ws.Range("D5").Value = "01/01"
$ws.Range("E5").Value = "01/07"

$ws.Range("D6").Value = "01/08"
$ws.Range("E6").Value = "01/14"

$ws.Range("D7").Value = "01/15"
$ws.Range("E7").Value = "01/21"

$ws.Range("D8").Value = "01/22"
$ws.Range("E8").Value = "01/28"

$ws.Range("D9").Value = "01/29"
$ws.Range("E9").Value = "02/04"

$ws.Range("D10").Value = "02/05"
$ws.Range("E10").Value = "02/11"

$ws.Range("D11").Value = "02/12"
$ws.Range("E11").Value = "02/18"

$ws.Range("D12").Value = "02/19"
$ws.Range("E12").Value = "02/25"

$ws.Range("D13").Value = "02/26"
$ws.Range("E13").Value = "03/04"
